$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A34").Value = 'Golang Architect / Principal Backend Architect - Atlanta, GA'
$ws.Range("B34").Value = 'https://www.dice.com/job-detail/9c4df3f1-bfaa-4b3b-ae7e-91d9ba3accec'
$ws.Range("C34").Value = 'Atlanta, Georgia'
$ws.Range("D34").Value = 'Contract, Third Party'
$ws.Range("E34").Value = '$80 - $85'
$ws.Range("F34").Value = 'Ocean Blue Solution'
